{"js": "// The new bullet goes right after the last paragraph in the document\n// (\"Researched and completed the change management plan portion of the\n// Project Plan.\"). insertParagraph(\"...\", \"After\") on that paragraph\n// creates a sibling that inherits its list formatting (ListParagraph\n// style, numId 1 / ilvl 0) and fills it with the new run text.\nconst body = context.document.body;\nconst lastParagraph = body.paragraphs.getLast();\nlastParagraph.insertParagraph(\n  \"Created a product change log and uploaded it to teams.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The new bullet goes right after the last paragraph in the document\n# (\"Researched and completed the change management plan portion of the\n# Project Plan.\"), inheriting that paragraph's list formatting\n# (ListParagraph style, numId 1 / ilvl 0).\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph.Range.Text = \"Created a product change log and uploaded it to teams.\"\n"}
